$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("qt7")

# ---------------------------------------------------------------------
# Row 10 / Row 11 get tidied up to a uniform centered/bordered look
# (borrow the existing "center, thin border, no fill" format already
# used elsewhere in the table, e.g. B2) while B10/E10 keep their own
# highlight fill untouched. New column G placeholders are added too.
# ---------------------------------------------------------------------
$fmtSrc = $ws.Range("B2")
$fmtSrc.Copy()
foreach ($addr in "A10", "C10", "D10", "F10", "A11", "B11", "C11", "D11", "E11", "F11") {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# New column G placeholders on rows 10/11 (blank).
$ws.Range("G10").Font.Name = "新細明體"
$ws.Range("G11").Font.Name = "新細明體"

# Row 10/11 no longer carry an explicit 14.4pt height -> autofit drops
# back to the sheet default (13.8pt), matching the target XML.
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------------
# Row 12: blank spacer row across A:G.
# ---------------------------------------------------------------------
foreach ($col in "A", "B", "C", "D", "E", "F", "G") {
    $ws.Range($col + "12").Font.Name = "新細明體"
}

# ---------------------------------------------------------------------
# Rows 13/14: new "Slider LED" t15/t9 wiring callback table.
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "Slider LED"
$ws.Range("B13").Value = "B7"
$ws.Range("C13").Value = "A2"
$ws.Range("D13").Value = "c3"
$ws.Range("E13").Value = "c2"
$ws.Range("F13").Value = "c1"
$ws.Range("G13").Value = "c0"

$ws.Range("A14").Value = "btn"
$ws.Range("B14").Value = "b4"
$ws.Range("C14").Value = "a1"
foreach ($col in "D", "E", "F", "G") {
    $ws.Range($col + "14").Font.Name = "新細明體"
}

$ws.Range("A13:G14").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# Row 15: blank spacer row across A:G.
# ---------------------------------------------------------------------
foreach ($col in "A", "B", "C", "D", "E", "F", "G") {
    $ws.Range($col + "15").Font.Name = "新細明體"
}

# ---------------------------------------------------------------------
# Selection / view state.
# ---------------------------------------------------------------------
$ws.Range("C8").Select()
